$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999510171489,
    0.99819163496484609,
    0.99262845713997794,
    0.9943847472776175,
    0.98294580930062359,
    0.95646547170402674,
    0.95346312494462571,
    0.94721208802886059,
    0.94161987292252425,
    0.93693778353721724,
    0.93612171086588447,
    0.93482551354025234,
    0.93248333530603911,
    0.93111280111730554,
    0.93118288385331449,
    0.9318524981496169,
    0.92814472053143815,
    0.92703580788949935,
    0.99712131922179914,
    0.98968991609489687,
    0.98845854029273728,
    0.9871940339171803,
    0.9855941948842375,
    0.97257386728184148,
    0.96611695906084827,
    0.9532787167137069,
    0.94908408553291967,
    0.93154576627129004,
    0.91917105001109056,
    0.9136947719951799,
    0.90993090687575473,
    0.90825158798147521,
    0.90773157254990977
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
